$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login Test")

# --- Add new test case row (row 7) below the existing LOG-TC-03 row (row 6) ---
# Mirrors how the row was authored: copy the previous row's formatting down,
# then fill in the new test-case values (module/scenario IDs continue the
# existing LOG-0x / LOG-TS-0x series).

# 1) Seed row 7 with row 6's formatting (borders/fills/fonts/number formats).
$ws.Range("B6:L6").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# 2) Fill in the new test case's content.
$ws.Range("B7").Value = "LOG-04"
$ws.Range("C7").Value = "Verify required fields validation"
$ws.Range("D7").Value = "LOG-TC-04"
$ws.Range("E7").Value = "Login Module"
$ws.Range("F7").Value = "LOG-TS-04"
$ws.Range("G7").Value = "Login with empty email"
$ws.Range("H7").Value = "On login page "
$ws.Range("I7").Value = "Leave email empty, enter password, click Login"
$ws.Range("J7").Value = "`"`"`nPassword124"
$ws.Range("K7").Value = "Validation message appears"
$ws.Range("L7").Value = "High"

# 3) Row height: let Excel settle on the same auto height used for this row.
$ws.Rows.Item(7).RowHeight = 31.5

# 4) The Test Data cell (J7) carries a mailto hyperlink, same as the rows
#    above it. Add it (this also stamps the cached display text), then
#    restore the cell's real text since the hyperlink's cached display
#    text is independent of the cell value.
$ws.Hyperlinks.Add($ws.Range("J7"), "mailto:unkown@gmail.com%0aPassword123", "", "", "unkown@gmail.com`nPassword123")
$ws.Range("J7").Value = "`"`"`nPassword124"

# 5) Re-apply row 6's formatting on top so every cell in the new row -
#    including the hyperlink cell - keeps the same styling as the rest
#    of the table (Hyperlinks.Add resets the cell style otherwise).
$ws.Range("B6:L6").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# --- Selection bookkeeping: the sheet was left with column D selected ---
$ws.Columns.Item(4).Select()
